# Update "想去人数" (want-to-go count) values in column F for specific
# event rows on the "展览" and "全部类型" sheets, per the gh-pages
# regeneration commit (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheetId 1) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value  = 991   # was 990
$wsExpo.Range("F6").Value  = 5503  # was 5498
$wsExpo.Range("F8").Value  = 689   # was 690
$wsExpo.Range("F17").Value = 1844  # was 1841
$wsExpo.Range("F19").Value = 913   # was 912
$wsExpo.Range("F23").Value = 549   # was 548
$wsExpo.Range("F24").Value = 152   # was 151
$wsExpo.Range("F28").Value = 2895  # was 2893
$wsExpo.Range("F29").Value = 179   # was 178
$wsExpo.Range("F30").Value = 102   # was 100
$wsExpo.Range("F40").Value = 720   # was 719
$wsExpo.Range("F41").Value = 89    # was 88

# --- Sheet "全部类型" (sheetId 4) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 991   # was 990
$wsAll.Range("F7").Value  = 5503  # was 5498
$wsAll.Range("F9").Value  = 689   # was 690
$wsAll.Range("F11").Value = 189   # was 188
$wsAll.Range("F23").Value = 1844  # was 1841
$wsAll.Range("F25").Value = 913   # was 912
$wsAll.Range("F29").Value = 549   # was 548
$wsAll.Range("F30").Value = 152   # was 151
$wsAll.Range("F32").Value = 2895  # was 2893
$wsAll.Range("F33").Value = 179   # was 178
$wsAll.Range("F34").Value = 102   # was 100
$wsAll.Range("F43").Value = 720   # was 719
$wsAll.Range("F44").Value = 89    # was 88

$wb.Save()
